$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

$row = 80

$ws.Cells.Item($row, 1).Value = "Sam Walton, Made in America"
$ws.Cells.Item($row, 2).Value = "Sam Walton"

$ws.Cells.Item($row - 1, 3).Copy()
$ws.Cells.Item($row, 3).PasteSpecial(-4122)
$ws.Cells.Item($row, 3).Value = "5/13/2020"

$ws.Cells.Item($row - 1, 4).Copy()
$ws.Cells.Item($row, 4).PasteSpecial(-4122)
$ws.Cells.Item($row, 4).Value = "5/27/2020"

$ws.Cells.Item($row, 5).Value = "biography;sam walton;business;walmart;retail"
$ws.Cells.Item($row, 6).Value = "Ebook"
$ws.Cells.Item($row, 7).Value = "379 Pages"

$ws.Application.CutCopyMode = $false

$ws.Range("A81").Select()
